$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 410
$ws.Range("I2").Value = 387.5
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 387.5
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -274.5
$ws.Range("N2").Value = -726
$ws.Range("H10").Value = 16661.75
$ws.Range("I10").Value = 6662.3335
$ws.Range("K10").Value = 6662.3335
$ws.Range("M10").Value = -6369.3335
$ws.Range("H18").Value = 1264.1428
$ws.Range("I18").Value = 1307.8334
$ws.Range("K18").Value = 1307.8334
$ws.Range("M18").Value = -1023.8334
$ws.Range("H51").Value = 4767335.5
$ws.Range("J51").Value = 7148956
$ws.Range("L51").Value = 7148956
$ws.Range("N51").Value = -7149924
$ws.Range("H58").Value = 2268.9375
$ws.Range("I58").Value = 423.30768
$ws.Range("J58").Value = 10266.667
$ws.Range("K58").Value = 1269.92304
$ws.Range("L58").Value = 30800.001
$ws.Range("M58").Value = -1119.92304
$ws.Range("N58").Value = -31100.001
$ws.Range("H100").Value = 1908.6154
$ws.Range("I100").Value = 1346.5454
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1346.5454
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -805.5454
$ws.Range("N100").Value = -6082
$ws.Range("H112").Value = 1972.3715
$ws.Range("J112").Value = 1824.5
$ws.Range("L112").Value = 5473.5
$ws.Range("N112").Value = -7689.5
$ws.Range("H129").Value = 704.63635
$ws.Range("J129").Value = 1697.6666
$ws.Range("L129").Value = 5092.9998
$ws.Range("N129").Value = -15092.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 6961.2666
$ws.Range("I4").Value = 289.33334
$ws.Range("J4").Value = 11409.223
$ws.Range("K4").Value = 289.33334
$ws.Range("L4").Value = 11409.223
$ws.Range("M4").Value = -173.33334
$ws.Range("N4").Value = -11641.223
$ws.Range("H11").Value = 9292001
$ws.Range("I11").Value = 13754751
$ws.Range("J11").Value = 3341668
$ws.Range("K11").Value = 13754751
$ws.Range("L11").Value = 3341668
$ws.Range("M11").Value = -13754607
$ws.Range("N11").Value = -3341956
$ws.Range("H97").Value = 1030.125
$ws.Range("I97").Value = 1044.8
$ws.Range("K97").Value = 1044.8
$ws.Range("M97").Value = -548.8
$ws.Range("H102").Value = 11693.667
$ws.Range("I102").Value = 3248.5
$ws.Range("J102").Value = 41251.75
$ws.Range("K102").Value = 3248.5
$ws.Range("L102").Value = 41251.75
$ws.Range("M102").Value = -1626.5
$ws.Range("N102").Value = -44495.75
$ws.Range("H122").Value = 19611380
$ws.Range("I122").Value = 33336096
$ws.Range("K122").Value = 100008288
$ws.Range("M122").Value = -100005838
$ws.Range("H132").Value = 166669500
$ws.Range("I132").Value = 200002800
$ws.Range("K132").Value = 600008400
$ws.Range("M132").Value = -600005870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 230.125
$ws.Range("I22").Value = 320.33334
$ws.Range("K22").Value = 320.33334
$ws.Range("M22").Value = -147.33334
$ws.Range("H23").Value = 1799
$ws.Range("J23").Value = 1799
$ws.Range("L23").Value = 1799
$ws.Range("N23").Value = -2365
$ws.Range("H105").Value = 1863.375
$ws.Range("I105").Value = 1819.1666
$ws.Range("J105").Value = 1996
$ws.Range("K105").Value = 1819.1666
$ws.Range("L105").Value = 1996
$ws.Range("M105").Value = -72.16660000000002
$ws.Range("N105").Value = -5490
$ws.Range("H107").Value = 13514305
$ws.Range("I107").Value = 850.6774
$ws.Range("K107").Value = 850.6774
$ws.Range("M107").Value = 1069.3226
$ws.Range("H134").Value = 1603
$ws.Range("I134").Value = 1440.875
$ws.Range("K134").Value = 4322.625
$ws.Range("M134").Value = -1787.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2948.875
$ws.Range("J5").Value = 3097.3333
$ws.Range("L5").Value = 3097.3333
$ws.Range("N5").Value = -3321.3333
$ws.Range("H6").Value = 3216223.8
$ws.Range("I6").Value = 4501914
$ws.Range("J6").Value = 1998.5
$ws.Range("K6").Value = 4501914
$ws.Range("L6").Value = 1998.5
$ws.Range("M6").Value = -4501801
$ws.Range("N6").Value = -2224.5
$ws.Range("H11").Value = 4645.778
$ws.Range("I11").Value = 3997
$ws.Range("K11").Value = 3997
$ws.Range("M11").Value = -3857
$ws.Range("H12").Value = 2696.6
$ws.Range("I12").Value = 2499.3333
$ws.Range("J12").Value = 2992.5
$ws.Range("K12").Value = 2499.3333
$ws.Range("L12").Value = 2992.5
$ws.Range("M12").Value = -2329.3333
$ws.Range("N12").Value = -3332.5
$ws.Range("H107").Value = 2132.8823
$ws.Range("I107").Value = 2473.1428
$ws.Range("J107").Value = 545
$ws.Range("K107").Value = 2473.1428
$ws.Range("L107").Value = 545
$ws.Range("M107").Value = -553.1428000000001
$ws.Range("N107").Value = -4385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17380484
$ws.Range("I4").Value = 42359056
$ws.Range("K4").Value = 127077168
$ws.Range("M4").Value = -127077056
$ws.Range("H8").Value = 27139.834
$ws.Range("I8").Value = 27139.834
$ws.Range("K8").Value = 81419.50199999999
$ws.Range("M8").Value = -81280.50199999999
$ws.Range("H86").Value = 1822.2222
$ws.Range("I86").Value = 667
$ws.Range("K86").Value = 2001
$ws.Range("M86").Value = -815
$ws.Range("H89").Value = 1822.2222
$ws.Range("I89").Value = 667
$ws.Range("K89").Value = 6003
$ws.Range("M89").Value = -75
$ws.Range("H132").Value = 2172.2856
$ws.Range("I132").Value = 1467
$ws.Range("J132").Value = 2701.25
$ws.Range("K132").Value = 13203
$ws.Range("L132").Value = 24311.25
$ws.Range("M132").Value = -10673
$ws.Range("N132").Value = -29371.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 636.1579
$ws.Range("J2").Value = 182.57143
$ws.Range("L2").Value = 182.57143
$ws.Range("N2").Value = -408.57143
$ws.Range("H68").Value = 100000
$ws.Range("I68").Value = 100000
$ws.Range("K68").Value = 100000
$ws.Range("M68").Value = -99189
$ws.Range("H71").Value = 100000
$ws.Range("I71").Value = 100000
$ws.Range("K71").Value = 300000
$ws.Range("M71").Value = -295944
$ws.Range("H97").Value = 1788.963
$ws.Range("I97").Value = 1607.55
$ws.Range("K97").Value = 1607.55
$ws.Range("M97").Value = -1111.55
$ws.Range("H102").Value = 1477.7778
$ws.Range("I102").Value = 722.7308
$ws.Range("K102").Value = 722.7308
$ws.Range("M102").Value = 899.2692
$ws.Range("H132").Value = 2798.186
$ws.Range("I132").Value = 2439.2222
$ws.Range("K132").Value = 7317.6666
$ws.Range("M132").Value = -4787.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1576.5
$ws.Range("I19").Value = 985
$ws.Range("J19").Value = 1872.25
$ws.Range("K19").Value = 985
$ws.Range("L19").Value = 1872.25
$ws.Range("M19").Value = -815
$ws.Range("N19").Value = -2212.25
$ws.Range("H40").Value = 2231.4
$ws.Range("I40").Value = 2034.8889
$ws.Range("K40").Value = 2034.8889
$ws.Range("M40").Value = -1898.8889
$ws.Range("H100").Value = 2134.037
$ws.Range("I100").Value = 1830.7059
$ws.Range("J100").Value = 2649.7
$ws.Range("K100").Value = 1830.7059
$ws.Range("L100").Value = 2649.7
$ws.Range("M100").Value = -1289.7059
$ws.Range("N100").Value = -3731.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1990.4706
$ws.Range("I122").Value = 1939.875
$ws.Range("K122").Value = 5819.625
$ws.Range("M122").Value = -3369.625
